# "Fruta / hortaliza, semanal" weekly update:
# a new weekly price record for Apio (Femacal de La Calera) is inserted
# at row 223, pushing the existing rows 223-261 down to 224-262.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 223 (shifts 223:261 -> 224:262).
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(223, 1).Value  = 3
$ws.Cells.Item(223, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(223, 3).Value  = "Coquimbo"
$ws.Cells.Item(223, 4).Value  = 44504
$ws.Cells.Item(223, 5).Value  = 5
$ws.Cells.Item(223, 6).Value  = 100112017
$ws.Cells.Item(223, 7).Value  = "Apio"
$ws.Cells.Item(223, 8).Value  = "Americana (o)"
$ws.Cells.Item(223, 9).Value  = "Primera"
$ws.Cells.Item(223, 10).Value = 160
$ws.Cells.Item(223, 11).Value = 9000
$ws.Cells.Item(223, 12).Value = 9000
$ws.Cells.Item(223, 13).Value = 9000
$ws.Cells.Item(223, 14).Value = "`$/docena de matas"
$ws.Cells.Item(223, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(223, 16).Value = 1500
$ws.Cells.Item(223, 17).Value = 6
$ws.Cells.Item(223, 18).Value = "Hortaliza"
